$d = $word.ActiveDocument

$replacements = @(
    @("49×51=2499", "59×25=1475"),
    @("60×40=2400", "98×66=6468"),
    @("96×32=3072", "57×51=2907"),
    @("35×60=2100", "27×96=2592"),
    @("46×89=4094", "30×17=510"),
    @("40×17=680",  "26×66=1716"),
    @("94×62=5828", "43×89=3827"),
    @("92×51=4692", "79×57=4503"),
    @("69×30=2070", "17×49=833"),
    @("73×92=6716", "70×83=5810"),
    @("20×41=820",  "81×53=4293"),
    @("40×48=1920", "82×31=2542"),
    @("56×76=4256", "25×14=350"),
    @("16×62=992",  "18×66=1188"),
    @("12×58=696",  "38×36=1368"),
    @("36×53=1908", "71×71=5041"),
    @("52×68=3536", "87×40=3480"),
    @("50×99=4950", "20×21=420"),
    @("59×37=2183", "27×28=756"),
    @("54×56=3024", "94×58=5452"),
    @("95×16=1520", "31×36=1116"),
    @("63×33=2079", "45×55=2475"),
    @("60×18=1080", "81×81=6561"),
    @("20×18=360",  "52×83=4316"),
    @("78×53=4134", "88×42=3696")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false,
                                  $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
